# This script fills in the exam-duty bill form (MsNaziaJahanKhanChowdhury.xlsx)
# with the teacher's name, designation, department, year/term, and the
# total-amount-in-words, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Name field (merged A3:C3) ---
$ws.Range("A3").Value = "নাম: Ms. Nazia Jahan Khan Chowdhury"

# --- Row 4: Designation field (merged A4:C4) + Year/Term values ---
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("G4").Value = "৪র্থ"
$ws.Range("I4").Value = "১ম"

# --- Row 5: Department short code + department field (merged F5:I5) ---
$ws.Range("B5").Value = "সিএসই"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# --- Row 32: Total amount in words (merged A32:E32) ---
$ws.Range("A32").Value = "কথায়:চৌদ্দ হাজার নয়শো চল্লিশ টাকা মাত্র।"

# --- Formatting tweaks ---
# Column A widened to fit the name/designation text
$ws.Columns.Item(1).ColumnWidth = 13.5

# Row 36 grown taller
$ws.Rows.Item(36).RowHeight = 68.4

# --- View state: select I32 (best effort for on-screen scroll position) ---
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("I32").Select()
